$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# cm007 (row 9) is being added/linked into the syllabus:
# - topic text changes from "Debugging errors" to "Debugging and defensive programming"
# - link_it flag flips from FALSE to TRUE
$ws.Range("D9").Value = "Debugging and defensive programming"
$ws.Range("C9").Value = $true

# Update the active selection to reflect the next cell (C10)
$ws.Range("C10").Select() | Out-Null
